# "Generate Report for Archive"
#
# The localization status flips from "Ready for handoff" to "In Translation"
# everywhere it is used (the Overview sheet's per-locale status columns, plus
# the Status column on each locale detail sheet). The shorter replacement
# text then lets the two affected status columns narrow on each sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: Status column (col C) ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: Status column (col C) ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Columns.Item(3).ColumnWidth = 12.5
